$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "05bd14b0-..." row Latest HO Xliff Generate Date
$wsOverview.Range("G3").Value = "2016-08-17 20:47:12"

# de-de sheet: same row/value ("Correspond Handoff Datetime") mirrors Overview's G3
$wsDeDe.Range("H3").Value = "2016-08-17 20:47:12"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback Datetime for row 3
$wsZhCn.Range("H3").Value = "2016-08-17 20:47:01"
$wsZhCn.Range("K3").Value = "2016-08-17 20:47:31"

# de-de sheet: Correspond Handback Datetime for row 3
$wsDeDe.Range("K3").Value = "2016-08-17 20:47:39"
